$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.840.13"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.296.08"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +18.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "269.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.01%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "2.637.16"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "2.282.93"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "43.702.06"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +13.53%  "
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  +7.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0935"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.127"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.43%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.106"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +16.50%  "
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.75%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.08%  "
